$d = $word.ActiveDocument

# Locate the paragraph that contains the "upload deadline" sentence so the
# three word-level substitutions below only ever touch that sentence, even
# if any of the search terms happen to appear elsewhere in the template.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*applicant and respondent may upload any written submissions*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the applicant/respondent upload-deadline paragraph"
}

# applicant -> claimant, respondent -> defendant, may -> should
# (mirrors the substantive wording change called out in the commit message)
$r1 = $target.Range
[void]$r1.Find.Execute("applicant", $true, $false, $false, $false, $false, $true, 1, $false, "claimant", 2)

$r2 = $target.Range
[void]$r2.Find.Execute("respondent", $true, $false, $false, $false, $false, $true, 1, $false, "defendant", 2)

$r3 = $target.Range
[void]$r3.Find.Execute("may upload", $true, $false, $false, $false, $false, $true, 1, $false, "should upload", 2)

Write-Host "Updated sentence: $($target.Range.Text)"
